$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 467, shifting rows 467:520 down to 468:521.
$ws.Rows("467:467").Insert()

# Fill in the new row 467 with its data.
$ws.Range("A467").Value = 4
$ws.Range("B467").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C467").Value = "Los Lagos"
$ws.Range("D467").Value = 45166
$ws.Range("E467").Value = 10
$ws.Range("F467").Value = "Fruta"
$ws.Range("G467").Value = 100102
$ws.Range("H467").Value = "Cítricos"
$ws.Range("I467").Value = 100102004
$ws.Range("J467").Value = "Mandarina"
$ws.Range("K467").Value = "Murcott"
$ws.Range("L467").Value = "Segunda"
$ws.Range("M467").Value = 250
$ws.Range("N467").Value = 8000
$ws.Range("O467").Value = 8000
$ws.Range("P467").Value = 8000
$ws.Range("Q467").Value = "$/bandeja 10 kilos"
$ws.Range("R467").Value = "Región de O'Higgins"
$ws.Range("S467").Value = 800
$ws.Range("T467").Value = 10
